# "Generate Report for Handback" - refresh the handback timestamps that
# were produced by the second (06c774b6...) pass of the report generator.
#
#   Overview!G3                       06:42:47 -> 06:43:30
#   zh-cn!H3  (Correspond Handoff)    06:42:42 -> 06:43:26
#   zh-cn!K3  (Correspond Handback)   06:43:00 -> 06:43:43
#   de-de!H3  (Correspond Handoff)    06:42:47 -> 06:43:30 (shares the
#                                      same shared-string as Overview!G3,
#                                      so it updates along with it)
#   de-de!K3  (Correspond Handback)   06:43:13 -> 06:43:50

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-23 06:43:30"

$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-23 06:43:26"
$wsZhCn.Range("K3").Value = "2016-08-23 06:43:43"

$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-23 06:43:30"
$wsDeDe.Range("K3").Value = "2016-08-23 06:43:50"
